$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 277.2857
$ws.Range("I42").Value = 51.75
$ws.Range("J42").Value = 367.5
$ws.Range("K42").Value = 155.25
$ws.Range("L42").Value = 1102.5
$ws.Range("M42").Value = 74.75
$ws.Range("N42").Value = -1562.5
$ws.Range("H80").Value = 1667.6666
$ws.Range("I80").Value = 1268.6666
$ws.Range("K80").Value = 3805.9998
$ws.Range("M80").Value = -2807.9998
$ws.Range("H83").Value = 1667.6666
$ws.Range("I83").Value = 1268.6666
$ws.Range("K83").Value = 11417.9994
$ws.Range("M83").Value = -6425.999400000001
$ws.Range("H88").Value = 1201.7778
$ws.Range("I88").Value = 1592.375
$ws.Range("J88").Value = 889.3
$ws.Range("K88").Value = 1592.375
$ws.Range("L88").Value = 889.3
$ws.Range("M88").Value = -1186.375
$ws.Range("N88").Value = -1701.3
$ws.Range("H91").Value = 1201.7778
$ws.Range("I91").Value = 1592.375
$ws.Range("J91").Value = 889.3
$ws.Range("K91").Value = 1592.375
$ws.Range("L91").Value = 889.3
$ws.Range("M91").Value = -188.375
$ws.Range("N91").Value = -3697.3
$ws.Range("H97").Value = 1440.5555
$ws.Range("I97").Value = 600
$ws.Range("K97").Value = 1800
$ws.Range("M97").Value = -1304
$ws.Range("H100").Value = 7617
$ws.Range("I100").Value = 7839.909
$ws.Range("J100").Value = 6799.6665
$ws.Range("K100").Value = 7839.909
$ws.Range("L100").Value = 6799.6665
$ws.Range("M100").Value = -7298.909
$ws.Range("N100").Value = -7881.6665
$ws.Range("H107").Value = 2617.2
$ws.Range("I107").Value = 1154.8
$ws.Range("K107").Value = 1154.8
$ws.Range("M107").Value = 765.2
$ws.Range("H113").Value = 116870.11
$ws.Range("I113").Value = 501666.5
$ws.Range("J113").Value = 6928.2856
$ws.Range("K113").Value = 501666.5
$ws.Range("L113").Value = 6928.2856
$ws.Range("M113").Value = -498412.5
$ws.Range("N113").Value = -13436.2856
$ws.Range("H137").Value = 1089.579
$ws.Range("I137").Value = 1013.25
$ws.Range("J137").Value = 1496.6666
$ws.Range("K137").Value = 3039.75
$ws.Range("L137").Value = 4489.9998
$ws.Range("M137").Value = -489.75
$ws.Range("N137").Value = -9589.9998
$ws.Range("H138").Value = 1988.4318
$ws.Range("I138").Value = 1271.2
$ws.Range("J138").Value = 2586.125
$ws.Range("K138").Value = 3813.6
$ws.Range("L138").Value = 7758.375
$ws.Range("M138").Value = 1326.4
$ws.Range("N138").Value = -18038.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5778.189
$ws.Range("I32").Value = 5296.879
$ws.Range("K32").Value = 5296.879
$ws.Range("M32").Value = -5009.879
$ws.Range("H45").Value = 12431
$ws.Range("I45").Value = 35670.332
$ws.Range("K45").Value = 35670.332
$ws.Range("M45").Value = -35293.332
$ws.Range("H61").Value = 2037.4
$ws.Range("I61").Value = 1422.8
$ws.Range("K61").Value = 1422.8
$ws.Range("M61").Value = -1210.8
$ws.Range("H64").Value = 12000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 12000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H110").Value = 1576.8462
$ws.Range("I110").Value = 1576.8462
$ws.Range("K110").Value = 1576.8462
$ws.Range("M110").Value = 468.1538
$ws.Range("H111").Value = 60322
$ws.Range("J111").Value = 60322
$ws.Range("L111").Value = 60322
$ws.Range("N111").Value = -68502
$ws.Range("H136").Value = 2037.4
$ws.Range("I136").Value = 1422.8
$ws.Range("K136").Value = 4268.4
$ws.Range("M136").Value = -1718.4
$ws.Range("H140").Value = 108411.2
$ws.Range("J140").Value = 108411.2
$ws.Range("L140").Value = 108411.2
$ws.Range("N140").Value = -118771.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7009.8184
$ws.Range("I20").Value = 8219.375
$ws.Range("J20").Value = 3784.3333
$ws.Range("K20").Value = 8219.375
$ws.Range("L20").Value = 3784.3333
$ws.Range("M20").Value = -7972.375
$ws.Range("N20").Value = -4278.3333
$ws.Range("H86").Value = 918.41174
$ws.Range("I86").Value = 915.2727
$ws.Range("J86").Value = 924.1667
$ws.Range("K86").Value = 915.2727
$ws.Range("L86").Value = 924.1667
$ws.Range("M86").Value = 207.7273
$ws.Range("N86").Value = -3170.1667
$ws.Range("H89").Value = 918.41174
$ws.Range("I89").Value = 915.2727
$ws.Range("J89").Value = 924.1667
$ws.Range("K89").Value = 4576.363499999999
$ws.Range("L89").Value = 4620.8335
$ws.Range("M89").Value = 1039.636500000001
$ws.Range("N89").Value = -15852.8335
$ws.Range("H134").Value = 1469.125
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null
$ws.Range("H141").Value = 67923.71000000001
$ws.Range("J141").Value = 67923.71000000001
$ws.Range("L141").Value = 67923.71000000001
$ws.Range("N141").Value = -78283.71000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13158.115
$ws.Range("I31").Value = 1481.3043
$ws.Range("K31").Value = 1481.3043
$ws.Range("M31").Value = -1186.3043
$ws.Range("H34").Value = 13158.115
$ws.Range("I34").Value = 1481.3043
$ws.Range("K34").Value = 1481.3043
$ws.Range("M34").Value = -1279.3043
$ws.Range("H86").Value = 4294.8
$ws.Range("I86").Value = 3831.3333
$ws.Range("K86").Value = 3831.3333
$ws.Range("M86").Value = -2708.3333
$ws.Range("H89").Value = 4294.8
$ws.Range("I89").Value = 3831.3333
$ws.Range("K89").Value = 19156.6665
$ws.Range("M89").Value = -13540.6665
$ws.Range("H141").Value = 335124.38
$ws.Range("J141").Value = 335124.38
$ws.Range("L141").Value = 335124.38
$ws.Range("N141").Value = -345484.38

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 5600
$ws.Range("I87").Value = 5600
$ws.Range("K87").Value = 16800
$ws.Range("M87").Value = -15552
$ws.Range("H90").Value = 5600
$ws.Range("I90").Value = 5600
$ws.Range("K90").Value = 50400
$ws.Range("M90").Value = -44160
$ws.Range("H113").Value = 1471.7587
$ws.Range("J113").Value = 1571.36
$ws.Range("L113").Value = 4714.08
$ws.Range("N113").Value = -9054.08
$ws.Range("H137").Value = 5558077.5
$ws.Range("I137").Value = 11112603
$ws.Range("K137").Value = 33337809
$ws.Range("M137").Value = -33332709

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7999
$ws.Range("I70").Value = 7999
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 7999
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -7729
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 7999
$ws.Range("I73").Value = 7999
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 7999
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -7063
$ws.Range("N73").Value = $null
$ws.Range("H107").Value = 125008900
$ws.Range("I107").Value = 300
$ws.Range("K107").Value = 300
$ws.Range("M107").Value = 1620
$ws.Range("H132").Value = 2116.2144
$ws.Range("I132").Value = 2116.2144
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6348.6432
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3818.6432
$ws.Range("N132").Value = $null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 850.5
$ws.Range("I22").Value = 801
$ws.Range("K22").Value = 801
$ws.Range("M22").Value = -506
$ws.Range("H27").Value = 850.5
$ws.Range("I27").Value = 801
$ws.Range("K27").Value = 801
$ws.Range("M27").Value = -694
$ws.Range("H40").Value = 40748.5
$ws.Range("I40").Value = 52499.5
$ws.Range("K40").Value = 52499.5
$ws.Range("M40").Value = -52363.5
$ws.Range("H46").Value = 47629.2
$ws.Range("I46").Value = 76115.336
$ws.Range("J46").Value = 4900
$ws.Range("K46").Value = 76115.336
$ws.Range("L46").Value = 4900
$ws.Range("M46").Value = -75927.336
$ws.Range("N46").Value = -5276
$ws.Range("H122").Value = 226888.78
$ws.Range("I122").Value = 336666.5
$ws.Range("K122").Value = 1009999.5
$ws.Range("M122").Value = -1007549.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 31747.5
$ws.Range("J49").Value = 31747.5
$ws.Range("L49").Value = 31747.5
$ws.Range("N49").Value = -32207.5
$ws.Range("H98").Value = 31090.5
$ws.Range("J98").Value = 31090.5
$ws.Range("L98").Value = 31090.5
$ws.Range("N98").Value = -37080.5
$ws.Range("H130").Value = 46224.668
$ws.Range("J130").Value = 46224.668
$ws.Range("L130").Value = 46224.668
$ws.Range("N130").Value = -56264.668
